$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 204, shifting existing rows 204-316 down to 205-317.
$ws.Rows.Item(204).Insert()

# Populate the newly inserted row 204 with the new record's data.
$ws.Range("A204").Value = 10
$ws.Range("B204").Value = "Vega Modelo de Temuco"
$ws.Range("C204").Value = "La Araucanía"
$ws.Range("D204").Value = 44606
$ws.Range("E204").Value = 9
$ws.Range("F204").Value = 100112037
$ws.Range("G204").Value = "Cebollín"
$ws.Range("H204").Value = "Sin especificar"
$ws.Range("I204").Value = "Primera"
$ws.Range("J204").Value = 90
$ws.Range("K204").Value = 8000
$ws.Range("L204").Value = 8000
$ws.Range("M204").Value = 8000
$ws.Range("N204").Value = "$/docena de paquetes"
$ws.Range("O204").Value = "Provincia de Cautín"
$ws.Range("P204").Value = 667
$ws.Range("Q204").Value = 12
$ws.Range("R204").Value = "Hortaliza"
